# Update the "取得日時" (acquired datetime) column on the "ランサーズ" sheet:
# all rows that currently show 2025-12-07 18:23:39 are refreshed to the
# latest scrape timestamp 2025-12-07 18:31:34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-07 18:23:39"
$newValue = "2025-12-07 18:31:34"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
